$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Range("B8").Value = "2026-01-01T13:37:23+00:00"
$ws1.Range("B21").Value = "'108"
$ws1.Range("B19").Copy()
$ws1.Range("B21").PasteSpecial(-4122)

$ws2 = $wb.Worksheets.Item("Concepts")

$ws2.Range("A54:D54").Copy($ws2.Range("A55:D55"))
$ws2.Cells.Item(55, 2).Value = "bristol-score"
$ws2.Cells.Item(55, 3).Value = "Bristol Stool Score"
$ws2.Cells.Item(55, 4).Value = "Bristol Stool Form Scale Score (1-7)"
$ws2.Range("A55:D55").Copy($ws2.Range("A56:D56"))
$ws2.Cells.Item(56, 2).Value = "bristol-1"
$ws2.Cells.Item(56, 3).Value = "Type 1"
$ws2.Cells.Item(56, 4).Value = "Separate hard lumps, like nuts (hard to pass)"
$ws2.Range("A56:D56").Copy($ws2.Range("A57:D57"))
$ws2.Cells.Item(57, 2).Value = "bristol-2"
$ws2.Cells.Item(57, 3).Value = "Type 2"
$ws2.Cells.Item(57, 4).Value = "Sausage-shaped but lumpy"
$ws2.Range("A57:D57").Copy($ws2.Range("A58:D58"))
$ws2.Cells.Item(58, 2).Value = "bristol-3"
$ws2.Cells.Item(58, 3).Value = "Type 3"
$ws2.Cells.Item(58, 4).Value = "Like a sausage but with cracks on its surface"
$ws2.Range("A58:D58").Copy($ws2.Range("A59:D59"))
$ws2.Cells.Item(59, 2).Value = "bristol-4"
$ws2.Cells.Item(59, 3).Value = "Type 4"
$ws2.Cells.Item(59, 4).Value = "Like a sausage or snake, smooth and soft"
$ws2.Range("A59:D59").Copy($ws2.Range("A60:D60"))
$ws2.Cells.Item(60, 2).Value = "bristol-5"
$ws2.Cells.Item(60, 3).Value = "Type 5"
$ws2.Cells.Item(60, 4).Value = "Soft blobs with clear-cut edges (passed easily)"
$ws2.Range("A60:D60").Copy($ws2.Range("A61:D61"))
$ws2.Cells.Item(61, 2).Value = "bristol-6"
$ws2.Cells.Item(61, 3).Value = "Type 6"
$ws2.Cells.Item(61, 4).Value = "Fluffy pieces with ragged edges, a mushy stool"
$ws2.Range("A61:D61").Copy($ws2.Range("A62:D62"))
$ws2.Cells.Item(62, 2).Value = "bristol-7"
$ws2.Cells.Item(62, 3).Value = "Type 7"
$ws2.Cells.Item(62, 4).Value = "Watery, no solid pieces. Entirely liquid"
$ws2.Range("A62:D62").Copy($ws2.Range("A63:D63"))
$ws2.Cells.Item(63, 2).Value = "abbey-score"
$ws2.Cells.Item(63, 3).Value = "Abbey Pain Scale Score"
$ws2.Cells.Item(63, 4).Value = "Total Abbey Pain Scale Score (0-100+ but usually 0-18+)"
$ws2.Range("A63:D63").Copy($ws2.Range("A64:D64"))
$ws2.Cells.Item(64, 2).Value = "abbey-vocalization"
$ws2.Cells.Item(64, 3).Value = "Vocalization"
$ws2.Cells.Item(64, 4).Value = "Whimpering, groaning, crying"
$ws2.Range("A64:D64").Copy($ws2.Range("A65:D65"))
$ws2.Cells.Item(65, 2).Value = "abbey-facial-expression"
$ws2.Cells.Item(65, 3).Value = "Facial Expression"
$ws2.Cells.Item(65, 4).Value = "Looking tense, frowning, grimacing, looking frightened"
$ws2.Range("A65:D65").Copy($ws2.Range("A66:D66"))
$ws2.Cells.Item(66, 2).Value = "abbey-body-language"
$ws2.Cells.Item(66, 3).Value = "Body Language"
$ws2.Cells.Item(66, 4).Value = "Fidgeting, rocking, guarding part of body, withdrawn"
$ws2.Range("A66:D66").Copy($ws2.Range("A67:D67"))
$ws2.Cells.Item(67, 2).Value = "abbey-behavioral-change"
$ws2.Cells.Item(67, 3).Value = "Behavioral Change"
$ws2.Cells.Item(67, 4).Value = "Increased confusion, refusing to eat, alteration in usual pattern"
$ws2.Range("A67:D67").Copy($ws2.Range("A68:D68"))
$ws2.Cells.Item(68, 2).Value = "abbey-psychological-change"
$ws2.Cells.Item(68, 3).Value = "Psychological Change"
$ws2.Cells.Item(68, 4).Value = "Temperature, pulse, blood pressure changes, perspiration, pallor"
$ws2.Range("A68:D68").Copy($ws2.Range("A69:D69"))
$ws2.Cells.Item(69, 2).Value = "abbey-physical-changes"
$ws2.Cells.Item(69, 3).Value = "Physical Changes"
$ws2.Cells.Item(69, 4).Value = "Skin tears, pressure areas, arthritis, contractures, previous injuries"
$ws2.Range("A69:D69").Copy($ws2.Range("A70:D70"))
$ws2.Cells.Item(70, 2).Value = "fluid-input-total"
$ws2.Cells.Item(70, 3).Value = "Total Fluid Input"
$ws2.Cells.Item(70, 4).Value = "Total fluid input over specified period (e.g. 24h)"
$ws2.Range("A70:D70").Copy($ws2.Range("A71:D71"))
$ws2.Cells.Item(71, 2).Value = "fluid-output-total"
$ws2.Cells.Item(71, 3).Value = "Total Fluid Output"
$ws2.Cells.Item(71, 4).Value = "Total fluid output over specified period (e.g. 24h)"
$ws2.Range("A71:D71").Copy($ws2.Range("A72:D72"))
$ws2.Cells.Item(72, 2).Value = "fluid-balance"
$ws2.Cells.Item(72, 3).Value = "Fluid Balance"
$ws2.Cells.Item(72, 4).Value = "Total Input minus Total Output"
$ws2.Range("A72:D72").Copy($ws2.Range("A73:D73"))
$ws2.Cells.Item(73, 2).Value = "urine-output"
$ws2.Cells.Item(73, 3).Value = "Urine Output"
$ws2.Cells.Item(73, 4).Value = "Volume of urine passed"
$ws2.Range("A73:D73").Copy($ws2.Range("A74:D74"))
$ws2.Cells.Item(74, 2).Value = "abc-chart"
$ws2.Cells.Item(74, 3).Value = "ABC Chart"
$ws2.Cells.Item(74, 4).Value = "Antecedent-Behaviour-Consequence Chart for PBS"
$ws2.Range("A74:D74").Copy($ws2.Range("A75:D75"))
$ws2.Cells.Item(75, 2).Value = "abc-antecedent"
$ws2.Cells.Item(75, 3).Value = "Antecedent"
$ws2.Cells.Item(75, 4).Value = "What happened immediately before the behaviour (triggers)"
$ws2.Range("A75:D75").Copy($ws2.Range("A76:D76"))
$ws2.Cells.Item(76, 2).Value = "abc-behaviour"
$ws2.Cells.Item(76, 3).Value = "Behaviour"
$ws2.Cells.Item(76, 4).Value = "Description of the behaviour itself (observable actions)"
$ws2.Range("A76:D76").Copy($ws2.Range("A77:D77"))
$ws2.Cells.Item(77, 2).Value = "abc-consequence"
$ws2.Cells.Item(77, 3).Value = "Consequence"
$ws2.Cells.Item(77, 4).Value = "What happened immediately after (response/outcome)"
$ws2.Range("A77:D77").Copy($ws2.Range("A78:D78"))
$ws2.Cells.Item(78, 2).Value = "abc-function"
$ws2.Cells.Item(78, 3).Value = "Function of Behaviour"
$ws2.Cells.Item(78, 4).Value = "Hypothesized function (e.g. Sensory, Escape, Attention, Tangible)"
$ws2.Range("A78:D78").Copy($ws2.Range("A79:D79"))
$ws2.Cells.Item(79, 2).Value = "abc-duration"
$ws2.Cells.Item(79, 3).Value = "Duration"
$ws2.Cells.Item(79, 4).Value = "Duration of the episode"
$ws2.Range("A79:D79").Copy($ws2.Range("A80:D80"))
$ws2.Cells.Item(80, 2).Value = "abc-intensity"
$ws2.Cells.Item(80, 3).Value = "Intensity"
$ws2.Cells.Item(80, 4).Value = "Intensity of the behaviour (1-10)"
$ws2.Range("A80:D80").Copy($ws2.Range("A81:D81"))
$ws2.Cells.Item(81, 2).Value = "oral-health-score"
$ws2.Cells.Item(81, 3).Value = "Oral Health Score"
$ws2.Cells.Item(81, 4).Value = "Total Oral Health Assessment Score"
$ws2.Range("A81:D81").Copy($ws2.Range("A82:D82"))
$ws2.Cells.Item(82, 2).Value = "oral-lips"
$ws2.Cells.Item(82, 3).Value = "Lips"
$ws2.Cells.Item(82, 4).Value = "Condition of lips (Pink/Moist vs Dry/Cracked)"
$ws2.Range("A82:D82").Copy($ws2.Range("A83:D83"))
$ws2.Cells.Item(83, 2).Value = "oral-tongue"
$ws2.Cells.Item(83, 3).Value = "Tongue"
$ws2.Cells.Item(83, 4).Value = "Condition of tongue (Pink/Moist vs Coated/Red)"
$ws2.Range("A83:D83").Copy($ws2.Range("A84:D84"))
$ws2.Cells.Item(84, 2).Value = "oral-gums"
$ws2.Cells.Item(84, 3).Value = "Gums"
$ws2.Cells.Item(84, 4).Value = "Condition of gums (Pink/Firm vs Bleeding/Receding)"
$ws2.Range("A84:D84").Copy($ws2.Range("A85:D85"))
$ws2.Cells.Item(85, 2).Value = "oral-teeth"
$ws2.Cells.Item(85, 3).Value = "Teeth/Dentures"
$ws2.Cells.Item(85, 4).Value = "Condition of teeth or dentures (Clean/Intact vs Decayed/Broken/Loose)"
$ws2.Range("A85:D85").Copy($ws2.Range("A86:D86"))
$ws2.Cells.Item(86, 2).Value = "oral-saliva"
$ws2.Cells.Item(86, 3).Value = "Saliva"
$ws2.Cells.Item(86, 4).Value = "Saliva quality (Moist/Watery vs Thick/Sticky/Absent)"
$ws2.Range("A86:D86").Copy($ws2.Range("A87:D87"))
$ws2.Cells.Item(87, 2).Value = "seizure-record"
$ws2.Cells.Item(87, 3).Value = "Seizure Record"
$ws2.Cells.Item(87, 4).Value = "Record of a seizure event"
$ws2.Range("A87:D87").Copy($ws2.Range("A88:D88"))
$ws2.Cells.Item(88, 2).Value = "seizure-type"
$ws2.Cells.Item(88, 3).Value = "Seizure Type"
$ws2.Cells.Item(88, 4).Value = "Type of seizure (Tonic-Clonic, Absence, Focal, etc)"
$ws2.Range("A88:D88").Copy($ws2.Range("A89:D89"))
$ws2.Cells.Item(89, 2).Value = "seizure-duration"
$ws2.Cells.Item(89, 3).Value = "Seizure Duration"
$ws2.Cells.Item(89, 4).Value = "Duration of the active seizure phase"
$ws2.Range("A89:D89").Copy($ws2.Range("A90:D90"))
$ws2.Cells.Item(90, 2).Value = "seizure-recovery"
$ws2.Cells.Item(90, 3).Value = "Recovery Phase"
$ws2.Cells.Item(90, 4).Value = "Duration/Description of post-ictal recovery"
$ws2.Range("A90:D90").Copy($ws2.Range("A91:D91"))
$ws2.Cells.Item(91, 2).Value = "seizure-trigger"
$ws2.Cells.Item(91, 3).Value = "Trigger"
$ws2.Cells.Item(91, 4).Value = "Suspected trigger for the seizure"
$ws2.Range("A91:D91").Copy($ws2.Range("A92:D92"))
$ws2.Cells.Item(92, 2).Value = "sleep-record"
$ws2.Cells.Item(92, 3).Value = "Sleep Record"
$ws2.Cells.Item(92, 4).Value = "Record of a sleep period"
$ws2.Range("A92:D92").Copy($ws2.Range("A93:D93"))
$ws2.Cells.Item(93, 2).Value = "sleep-quality"
$ws2.Cells.Item(93, 3).Value = "Sleep Quality"
$ws2.Cells.Item(93, 4).Value = "Subjective or observed quality of sleep"
$ws2.Range("A93:D93").Copy($ws2.Range("A94:D94"))
$ws2.Cells.Item(94, 2).Value = "sleep-hours"
$ws2.Cells.Item(94, 3).Value = "Hours Slept"
$ws2.Cells.Item(94, 4).Value = "Total hours of sleep achieved"
$ws2.Range("A94:D94").Copy($ws2.Range("A95:D95"))
$ws2.Cells.Item(95, 2).Value = "sleep-disturbances"
$ws2.Cells.Item(95, 3).Value = "Disturbances"
$ws2.Cells.Item(95, 4).Value = "Number or description of distinct awakenings"
$ws2.Range("A95:D95").Copy($ws2.Range("A96:D96"))
$ws2.Cells.Item(96, 2).Value = "urinalysis-panel"
$ws2.Cells.Item(96, 3).Value = "Urinalysis Panel"
$ws2.Cells.Item(96, 4).Value = "Urine Dipstick Test Panel"
$ws2.Range("A96:D96").Copy($ws2.Range("A97:D97"))
$ws2.Cells.Item(97, 2).Value = "ua-leukocytes"
$ws2.Cells.Item(97, 3).Value = "Leukocytes"
$ws2.Cells.Item(97, 4).Value = "Leukocytes (WBCs) in urine"
$ws2.Range("A97:D97").Copy($ws2.Range("A98:D98"))
$ws2.Cells.Item(98, 2).Value = "ua-nitrites"
$ws2.Cells.Item(98, 3).Value = "Nitrites"
$ws2.Cells.Item(98, 4).Value = "Nitrites in urine"
$ws2.Range("A98:D98").Copy($ws2.Range("A99:D99"))
$ws2.Cells.Item(99, 2).Value = "ua-protein"
$ws2.Cells.Item(99, 3).Value = "Protein"
$ws2.Cells.Item(99, 4).Value = "Protein in urine"
$ws2.Range("A99:D99").Copy($ws2.Range("A100:D100"))
$ws2.Cells.Item(100, 2).Value = "ua-blood"
$ws2.Cells.Item(100, 3).Value = "Blood"
$ws2.Cells.Item(100, 4).Value = "Blood (Hemoglobin) in urine"
$ws2.Range("A100:D100").Copy($ws2.Range("A101:D101"))
$ws2.Cells.Item(101, 2).Value = "ua-glucose"
$ws2.Cells.Item(101, 3).Value = "Glucose"
$ws2.Cells.Item(101, 4).Value = "Glucose in urine"
$ws2.Range("A101:D101").Copy($ws2.Range("A102:D102"))
$ws2.Cells.Item(102, 2).Value = "ua-ketones"
$ws2.Cells.Item(102, 3).Value = "Ketones"
$ws2.Cells.Item(102, 4).Value = "Ketones in urine"
$ws2.Range("A102:D102").Copy($ws2.Range("A103:D103"))
$ws2.Cells.Item(103, 2).Value = "ua-ph"
$ws2.Cells.Item(103, 3).Value = "pH"
$ws2.Cells.Item(103, 4).Value = "Urine pH Level"
$ws2.Range("A103:D103").Copy($ws2.Range("A104:D104"))
$ws2.Cells.Item(104, 2).Value = "ua-sg"
$ws2.Cells.Item(104, 3).Value = "Specific Gravity"
$ws2.Cells.Item(104, 4).Value = "Urine Specific Gravity"
$ws2.Range("A104:D104").Copy($ws2.Range("A105:D105"))
$ws2.Cells.Item(105, 2).Value = "mca-present"
$ws2.Cells.Item(105, 3).Value = "Capacity Present"
$ws2.Cells.Item(105, 4).Value = "Patient has capacity for this decision"
$ws2.Range("A105:D105").Copy($ws2.Range("A106:D106"))
$ws2.Cells.Item(106, 2).Value = "4at-change-no"
$ws2.Cells.Item(106, 3).Value = "No Acute Change"
$ws2.Cells.Item(106, 4).Value = "No indication of acute change or fluctuating course"
$ws2.Range("A106:D106").Copy($ws2.Range("A107:D107"))
$ws2.Cells.Item(107, 2).Value = "4at-amt4-1error"
$ws2.Cells.Item(107, 3).Value = "1 Error"
$ws2.Cells.Item(107, 4).Value = "1 error in AMT4 test"
$ws2.Range("A107:D107").Copy($ws2.Range("A108:D108"))
$ws2.Cells.Item(108, 2).Value = "4at-attention-gt7"
$ws2.Cells.Item(108, 3).Value = "Months Backwards < 7 months"
$ws2.Cells.Item(108, 4).Value = "Less than 7 months correctly recited backwards"
$ws2.Range("A108:D108").Copy($ws2.Range("A109:D109"))
$ws2.Cells.Item(109, 2).Value = "4at-alert-normal"
$ws2.Cells.Item(109, 3).Value = "Normal Alertness"
$ws2.Cells.Item(109, 4).Value = "Patient is fully alert"
